# Daniels assesments merged in to my structure
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Peer  and self assessment")

# Criterion 1 Online collaboration - Daniel (row 7)
$ws.Range("B7").Value = "Good"
$ws.Range("C7").Value = "Decent response times, however sometimes to long a response. Very active on discord."

# Criterion 1 International Collaboration - Daniel (row 20)
$ws.Range("B20").Value = "Excellent"
$ws.Range("C20").Value = "Active collaborator, motivated"

# Reflect the final active selection/scroll position like the saved file
$ws.Range("C20").Select()
